# Add variable indexing for block cloning:
# Insert a new paragraph "${blockVariable}" right before the "${/CLONEME}"
# paragraph, with a "_GoBack" bookmark placed right before the closing
# "}" (splitting the text into two runs). Adding the bookmark with this
# reserved name automatically relocates it away from its previous
# location (the "This should be deleted." paragraph), which removes the
# old bookmarkStart/bookmarkEnd pair from there.

$d = $word.ActiveDocument

# Locate the "${/CLONEME}" paragraph so we can insert a new paragraph right before it.
$closeCloneMe = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*`${/CLONEME}*") {
        $closeCloneMe = $p
        break
    }
}

$insertionPoint = $d.Range($closeCloneMe.Range.Start, $closeCloneMe.Range.Start)
$insertionPoint.InsertParagraphBefore()

# Re-locate the freshly inserted (now empty) paragraph - it sits right
# before "${/CLONEME}".
$closeCloneMe2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*`${/CLONEME}*") {
        $closeCloneMe2 = $p
        break
    }
}
$newParaIndex = $closeCloneMe2.Index - 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newRange = $newPara.Range
$newRange.Text = [char]36 + "{blockVariable}"

# Re-fetch the paragraph/range (text assignment can shift range ends).
$newPara2 = $d.Paragraphs.Item($newParaIndex)
$fullRange = $newPara2.Range
$fullText = $fullRange.Text
$closeBracePos = $fullText.IndexOf("}")
$bookmarkPos = $fullRange.Start + $closeBracePos
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
